$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the text of existing cell D5 (shared string content edit)
$ws.Range("D5").Value = "Design of implementation new sync objects "

# Add the three new rows
$ws.Range("A5").Copy()
$ws.Range("A6:A8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A6").Value = 41430
$ws.Range("B6").Value = 2
$ws.Range("D6").Value = "Implementation of mutexes"

$ws.Range("A7").Value = 41431
$ws.Range("B7").Value = 2
$ws.Range("D7").Value = "Implementation of mutexes"

$ws.Range("A8").Value = 41432
$ws.Range("B8").Value = 4
$ws.Range("D8").Value = "Update Manual"

# Update selection to B6
$ws.Range("B6").Select()
